# Update the division-problem worksheet table with newly generated values.
# The table is a 20-row x 5-column grid; the problem rows are rows 1, 5, 9,
# 13 and 17 (1-based), each holding 5 expressions of the form "NN÷N=".
# We target each cell directly (rather than a blanket Find/Replace) because
# some of the new values coincide with *other* old values elsewhere in the
# table (e.g. "79÷4=" becomes "85÷3=", which is also the original text of
# another cell), so a global replace-all could clobber the wrong cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @(1, 1, "13÷5=", "56÷2="),
    @(1, 2, "85÷3=", "42÷9="),
    @(1, 3, "13÷7=", "62÷5="),
    @(1, 4, "23÷2=", "47÷3="),
    @(1, 5, "27÷9=", "62÷9="),

    @(5, 1, "45÷6=", "92÷7="),
    @(5, 2, "49÷5=", "71÷3="),
    @(5, 3, "26÷9=", "30÷3="),
    @(5, 4, "11÷8=", "30÷4="),
    @(5, 5, "72÷7=", "16÷7="),

    @(9, 1, "70÷9=", "49÷8="),
    @(9, 2, "79÷4=", "85÷3="),
    @(9, 3, "37÷8=", "78÷8="),
    @(9, 4, "29÷7=", "27÷7="),
    @(9, 5, "88÷6=", "69÷2="),

    @(13, 1, "69÷6=", "49÷9="),
    @(13, 2, "75÷5=", "27÷3="),
    @(13, 3, "80÷9=", "38÷5="),
    @(13, 4, "19÷4=", "72÷8="),
    @(13, 5, "40÷4=", "14÷2="),

    @(17, 1, "50÷8=", "22÷4="),
    @(17, 2, "28÷3=", "38÷4="),
    @(17, 3, "43÷6=", "91÷8="),
    @(17, 4, "70÷6=", "97÷7="),
    @(17, 5, "53÷5=", "46÷2=")
)

foreach ($item in $replacements) {
    $row = $item[0]
    $col = $item[1]
    $new = $item[3]

    $cell = $t.Cell($row, $col)
    $r = $cell.Range
    # Trim the trailing cell-end/paragraph marks so we only replace the
    # visible text, preserving the run's formatting (font, size, etc).
    $r.End = $r.End - 1
    $r.Text = $new
}
